# Term 2.0.0 update: bump metadata, replace existing FSIII concepts with
# four (plus one) newly-minted UUID concepts, and re-add the original
# J1..J5 concepts as five freshly appended "Include from FSIII N" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------
# 2. Replace the "Value" concept code on the five existing
#    "Include from FSIII" sheets with the new UUID-based concept codes.
# ---------------------------------------------------------------------
$existingSheetNames = @(
    "Include from FSIII",
    "Include from FSIII 2",
    "Include from FSIII 3",
    "Include from FSIII 4",
    "Include from FSIII 5"
)
$newConceptCodes = @(
    "43c2b7f0-5e55-4627-8fcf-bdaf5a9d84ac",
    "1c850a09-aa49-4fae-9354-f932f13e030b",
    "462f9352-0129-4d8e-8c75-a6dfed78ddcf",
    "4571f168-a92a-4caf-8dc8-35f45c2a1cb4",
    "86b53158-6d05-412e-ad55-2e1fa26359b3"
)

for ($i = 0; $i -lt $existingSheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($existingSheetNames[$i])
    $ws.Range("C2").Value = $newConceptCodes[$i]
}

# ---------------------------------------------------------------------
# 3. Append five new sheets ("Include from FSIII 6" .. "10") carrying
#    forward the original J1..J5 concept codes, using the same
#    Property/Operation/Value + System URI layout as the existing ones.
# ---------------------------------------------------------------------
$oldConceptCodes = @("J1", "J2", "J3", "J4", "J5")

$lastExisting = $wb.Worksheets.Item("Include from FSIII 5")
for ($i = 0; $i -lt $oldConceptCodes.Length; $i++) {
    $lastExisting.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = "Include from FSIII " + (6 + $i)
    $newSheet.Range("C2").Value = $oldConceptCodes[$i]
    $lastExisting = $newSheet
}
